$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4-7 down to 5-8
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new record
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(4, 4).Value = 45225
$ws.Cells.Item(4, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = 100112017
$ws.Cells.Item(4, 7).Value = "Corazón de apio"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 60
$ws.Cells.Item(4, 11).Value = 1500
$ws.Cells.Item(4, 12).Value = 2000
$ws.Cells.Item(4, 13).Value = 1750
$ws.Cells.Item(4, 14).Value = "$/docena de matas"
$ws.Cells.Item(4, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(4, 16).Value = 292
$ws.Cells.Item(4, 17).Value = 6
$ws.Cells.Item(4, 18).Value = "Hortaliza"
